$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = "AutomatedTest/01"
$ws.Range("D2").Value = 45480
$ws.Range("E2").Value = "PO/AutomatedTest/2024/001"
$ws.Range("F2").Value = "AutomatedTest/F4748"
$ws.Range("H2").Value = "A01123G026"
$ws.Range("K2").Value = "KG"
$ws.Range("M2").Value = "GUNA BUMI UTAMA"
$ws.Range("N2").Value = "Z08-221110"

# --- Row 3 ---
$ws.Range("C3").Value = "AutomatedTest/02"
$ws.Range("D3").Value = 45480
$ws.Range("E3").Value = "PO/AutomatedTest/2024/002"
$ws.Range("F3").Value = "AutomatedTest/F4746"
$ws.Range("H3").Value = "A01123G027"
$ws.Range("K3").Value = "KG"
$ws.Range("M3").Value = "PERTAMINA PATRA NIAGA"
$ws.Range("N3").Value = "Z08-221221"

# --- Row 4 ---
$ws.Range("C4").Value = "AutomatedTest/03"
$ws.Range("D4").Value = 45480
$ws.Range("E4").Value = "PO/AutomatedTest/2024/003"
$ws.Range("F4").Value = "AutomatedTest/F5266"
$ws.Range("H4").Value = "A01123G025"
$ws.Range("K4").Value = "KG"
$ws.Range("M4").Value = "PRASADHA PAMUNAH LIMBAH INDUSTRI"
$ws.Range("N4").Value = "Z08-221112"

# --- Row 5 ---
$ws.Range("B5").Value = "KAEF"
$ws.Range("C5").Value = "AutomatedTest/04"
$ws.Range("D5").Value = 45480
$ws.Range("E5").Value = "PO/AutomatedTest/2024/004"
$ws.Range("F5").Value = "AutomatedTest/F4740"
$ws.Range("H5").Value = "A01123G024"
$ws.Range("K5").Value = "KG"
$ws.Range("M5").Value = "SATYA SAMITRA NIAGATAMA, PT"
$ws.Range("N5").Value = "Z08-221113"

# --- Formatting: E:F columns rows 2-5 get wrap/top alignment (style idx 6) ---
$wrapRange = $ws.Range("E2:F5")
$wrapRange.WrapText = $true
$wrapRange.VerticalAlignment = -4160

# --- Formatting: M2:M4 get left/top/wrap + number format (style idx 7) ---
$mRange = $ws.Range("M2:M4")
$mRange.WrapText = $true
$mRange.VerticalAlignment = -4160
$mRange.HorizontalAlignment = -4131
$mRange.NumberFormat = "#,##0.00"

# --- New blank rows 6 and 7 with F column styled like E:F above ---
$ws.Range("F6").WrapText = $true
$ws.Range("F6").VerticalAlignment = -4160
$ws.Range("F7").WrapText = $true
$ws.Range("F7").VerticalAlignment = -4160

# --- Row heights for rows 2-5 (30pt, matching wrapped two-line text) ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30

# --- Selection ---
$ws.Range("M9").Select()
